$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Rewrite the "In summary: ..." paragraph (5th paragraph) in place
#    via three small, unambiguous Find/Replace calls instead of one
#    giant literal, so each call has unique surrounding context.
# ---------------------------------------------------------------------

# 1a) "aynu-game-struct/thing" -> "aynu-game-struct/thing (such as Sylene)"
$null = $d.Content.Find.Execute(
    "aynu-game-struct/thing with abstract aynu game-value",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "aynu-game-struct/thing (such as Sylene) with abstract aynu game-value",
    2)

# 1b) insert an extra bracket-group right before the existing
#     "[Elysion][Tyrion][Sirion][Esoteria][Deklein]" group
$null = $d.Content.Find.Execute(
    "[Asakai] , [Elysion][Tyrion][Sirion][Esoteria][Deklein] - [more aynu code to develop] and more",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[Asakai] , [Elysion][Tyrion][Coronis][Artorius][Antorus][Sirion][Exodius][Azrafel], [Elysion][Tyrion][Sirion][Esoteria][Deklein] - [more aynu code to develop] and more",
    2)

# 1c) extend the tail of the paragraph with the new "and motivates
#     players to ..." clause plus the new trailing sentence
$newTail = "through aynu-code by collecting and possessing these objects] and motivates players to collect and possess them for either their abstract-aynu-game-value, [Elysion][Firdaws][Tyrion][Ero][Coronis], [aynu-code expressing ultimate game-value-things], game-play use, [abstract-aynu-code expressing ultimate desires that make the player happy in a way that transcends mortal comprehension and can only be encoded using aynu-code], [Artorius][Coronis][Ho][Megas][Coronis][Tyrion], or any other [abstract-aynu-theory-reasons], …, [more to develop] that results from and is created by the abstract aynu-code each aynu-game-struct possesses. There are many different kinds of [abstract-aynu-structs] that will appeal to different kinds of players for many types of reasons."
$null = $d.Content.Find.Execute(
    "through aynu-code] by collecting and possessing these objects.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $newTail,
    2)

# ---------------------------------------------------------------------
# 2) Insert two brand-new paragraphs right after the paragraph we just
#    edited: a blank spacer paragraph, then a paragraph of new body
#    text -- both ahead of the pre-existing blank paragraph that used
#    to directly follow the "In summary: ..." paragraph.
# ---------------------------------------------------------------------
$spacerTarget = $d.Paragraphs(6)
$spacerTarget.Range.InsertParagraphBefore()
$spacerTarget.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs(7)
$newPara.Range.InsertAfter("Players will want to collect many different kinds of abstract aynu-game-structs, and each player will have their own reasons for wanting a particular game-struct. Different aynu-code will create different reasons for wanting to collect an objects and while some aynu-structs will not be desired by everyone, almost every aynu-struct will have a player with a set of [aynu-theory-reasons/desires] that will desire the struct for those [abstract-aynu-theory-reasons/desires]. There are also many different game-play reasons that can be relevant when determining whether a struct is desirable.")
